$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.288738369941711
$ws.Range("B1").Value = 2.161528825759888
$ws.Range("C1").Value = 5.23831844329834
$ws.Range("D1").Value = 0.5777249336242676
$ws.Range("E1").Value = 0.6801949739456177
